# Estilizando a pagina de inicio
# Adds a new "Relatar bugs" feature row (row 12) to the cronograma sheet,
# and tweaks column widths / fill formatting to match the refreshed layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new row at position 12 (pushes old rows 12-16 down to 13-17) ---
$ws.Rows.Item(12).Insert()

# Copy the formatting from row 8 (B8:F8), which already has the
# "s=1 customFormat, ht=72, wrap/center/border" look we want for the new row,
# onto the freshly inserted row 12 cells.
$fmtSrc = $ws.Range("B8:F8")
$fmtDst = $ws.Range("B12:F12")
$fmtSrc.Copy()
$fmtDst.PasteSpecial(-4122)   # xlPasteFormats
$ws.Application.CutCopyMode = $false

# Row height to match the new (taller) content block.
$ws.Rows.Item(12).RowHeight = 72

# --- 2. Fill in the new row's content ---
$ws.Range("B12").Value = "Relatar bugs"
$ws.Range("C12").Value = "Todos"
$ws.Range("D12").Value = "Um botão que ficará no dropdown do menu do usuário, e ao clicar nele, será redirecionado a um forms do google que permitirá a enviar ou relatar algum bug encontrado no sistema"
$ws.Range("E12").Value = ""
$ws.Range("F12").Value = "Na 1º Versão"

# --- 3. Remove the grey fill from the lower block of rows (old rows 12-16,
#         now shifted to rows 13-17) so they match rows 2-11's plain/bordered look ---
$noFillRange1 = $ws.Range("B13:D17")
$noFillRange1.Interior.Pattern = -4142   # xlPatternNone
$noFillRange2 = $ws.Range("F13:F17")
$noFillRange2.Interior.Pattern = -4142   # xlPatternNone

# --- 4. Column layout: column E (5) becomes narrower than the other columns ---
$ws.Columns.Item(5).ColumnWidth = 19.5

# --- 5. Update the view: scroll so row 7 is at the top, and select B9 ---
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B9").Select()
